$d = $word.ActiveDocument

# The trailing "_GoBack" bookmark currently sits at the start of the
# last paragraph (right before the "<img />" runs). In the edited
# document it is relocated to the end of the brand-new "26) ..." paragraph
# that gets appended after it, so remove it from its current spot first.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# Append a new paragraph at the very end of the document body (after the
# "<img />" paragraph) describing step 26 (installing bootstrap), ending
# with the relocated "_GoBack" bookmark.
$endPos = $d.Content.End
$insertionPoint = $d.Range($endPos, $endPos)

$newParagraphXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:rPr><w:b/><w:sz w:val="32"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:b/><w:sz w:val="32"/></w:rPr><w:t xml:space="preserve">26) </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:rPr><w:b/><w:sz w:val="32"/></w:rPr><w:t>install</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:rPr><w:b/><w:sz w:val="32"/></w:rPr><w:t xml:space="preserve"> bootstrap in react app: </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>npm</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve"> install react-bootstrap bootstrap</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '</w:p>'

$insertionPoint.InsertXML($newParagraphXml)
